$d = $word.ActiveDocument

# The document ends with the "Tri des status par leur priorite ..." bullet
# paragraph, right before the sectPr. Append two new list items (same
# "Paragraphedeliste" style / numbering as the rest of the list) after it.

$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()
$d.Paragraphs.Last.Range.Text = "Commit change"

$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()
$d.Paragraphs.Last.Range.Text = "Debut d'association de of core et of front. Je dois lancer elasticsearch et of core (collecte) depuis of front, des le lancement du serveur tomcat."
